$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with a new "Approval" column ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F21"))
$col6 = $lo.ListColumns.Item(6)
$col6.Range.Cells.Item(1, 1).Value = "Approval"

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = "Approved"
}

# --- Column width for new column F ---
$ws.Columns.Item(6).ColumnWidth = 9.33
$ws.Columns.Item(6).HorizontalAlignment = -4108

# --- Re-assert the autofilter across the full (now 6-column) range ---
$ws.Range("A1:F21").AutoFilter() | Out-Null

$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=data!`$F`$1:`$F`$21")
$nm.Visible = $false

# --- Thin borders around every cell in the table (header + body, A:F) ---
$table = $ws.Range("A1:F21")
$table.Borders.LineStyle = 1

# --- Header row A1:E1 stays with plain thin-bordered style (no fill) ---

# --- Header cell F1: bold white text, blue fill, centered ---
$hdr = $ws.Range("F1")
$hdr.Font.Bold = $true
$hdr.Font.ThemeColor = 1
$hdr.Interior.ThemeColor = 5
$hdr.HorizontalAlignment = -4108

# --- Data cells F2:F21: light accent fill ---
$body = $ws.Range("F2:F21")
$body.Interior.ThemeColor = 8
$body.Interior.TintAndShade = 0.4

$ws.Range("K11").Select() | Out-Null
